$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "test_checking[chrome]"
$ws.Range("D2").Value = "FAILED"
$ws.Range("E2").Value = 0.03103728900896385
$ws.Range("F2").Value = "2022-08-19T14:44:56"
$ws.Range("G2").Value = "AssertionError: Please enter a valid URL... http://127.0.0.1:5500/sample.html is not a valid URL."
